$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7642615437507629
$ws.Range("B1").Value = 0.5908836722373962
$ws.Range("C1").Value = 1.230963706970215
$ws.Range("D1").Value = 3.706138372421265
$ws.Range("E1").Value = 1.516143798828125
